$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so that
# numeric-looking values (e.g. "298.82") are not auto-converted to numbers
# and do not lose precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.185.79'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '2.267.29'
$ws.Range("E3").Value = '  -2.42%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '298.82'
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("D6").Value = '95.51'
$ws.Range("E6").Value = '  -4.88%  '
$ws.Range("D7").Value = '0.495'
$ws.Range("E7").Value = '  -2.49%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").Value = '33.33'
$ws.Range("E10").Value = '  -3.23%  '
$ws.Range("D11").Value = '0.0789'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = '48.09'
$ws.Range("E12").Value = '  -7.76%  '
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '6.64'
$ws.Range("E14").Value = '  -1.52%  '
$ws.Range("D15").Value = '2.622.23'
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("D16").Value = '15.53'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").Value = '2.282.60'
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").Value = '0.781'
$ws.Range("E18").Value = '  -4.55%  '
$ws.Range("D19").Value = '42.111.19'
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("D20").Value = '11.69'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '0.0₃0889'
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").Value = '5.97'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D23").Value = '66.54'
$ws.Range("E23").Value = '  -3.82%  '
$ws.Range("D24").Value = '234.40'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Value = '1.95'
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").Value = '23.91'
$ws.Range("E28").Value = '  -5.55%  '
$ws.Range("E29").Value = '  +2.86%  '
$ws.Range("D30").Value = '168.18'
$ws.Range("E30").Value = '  +4.73%  '
$ws.Range("D31").Value = '33.80'
$ws.Range("E31").Value = '  -2.75%  '
$ws.Range("D32").Value = '9.16'
$ws.Range("E32").Value = '  -0.58%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '4.90'
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("D35").Value = '4.47'
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("D36").Value = '16.48'
$ws.Range("E36").Value = '  -3.05%  '
$ws.Range("E37").Value = '  -4.87%  '
$ws.Range("D38").Value = '0.0686'
$ws.Range("E38").Value = '  -4.52%  '
$ws.Range("D39").Value = '2.77'
$ws.Range("E39").Value = '  -3.70%  '
$ws.Range("D40").Value = '0.0985'
$ws.Range("E40").Value = '  -2.58%  '
$ws.Range("E41").Value = '  -2.59%  '
$ws.Range("D42").Value = '1.72'
$ws.Range("E42").Value = '  -5.50%  '
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  -5.09%  '
$ws.Range("D44").Value = '1.958.49'
$ws.Range("E44").Value = '  -2.24%  '
$ws.Range("D45").Value = '0.0277'
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("D46").Value = '17.40'
$ws.Range("E46").Value = '  -6.75%  '
$ws.Range("D47").Value = '9.54'
$ws.Range("E47").Value = '  -6.15%  '
$ws.Range("D48").Value = '2.77'
$ws.Range("E48").Value = '  -3.61%  '
$ws.Range("D49").Value = '2.493.90'
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("D50").Value = '52.21'
$ws.Range("E50").Value = '  -5.54%  '
$ws.Range("D51").Value = '4.53'
$ws.Range("E51").Value = '  -2.92%  '
